$d = $word.ActiveDocument

# 1) "Il programma non necessita di nessuna installazione. "
#    -> "Il programma non necessita di alcuna installazione:"
$d.Content.Find.Execute(
    "Il programma non necessita di nessuna installazione. ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Il programma non necessita di alcuna installazione:", 2) | Out-Null

# 2) "Basta digitare" -> "Basterà digitare"  (Basta -> Basterà)
$d.Content.Find.Execute(
    "Basta digitare", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Basterà digitare", 2) | Out-Null

# 3) " basta fare un " -> ", basterà fare "
$d.Content.Find.Execute(
    " basta fare un ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    ", basterà fare ", 2) | Out-Null

# 4) "ed eseguirà in automatico il comando sopra citato." -> "e il programma verrà avviato."
$d.Content.Find.Execute(
    "ed eseguirà in automatico il comando sopra citato.", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "e il programma verrà avviato.", 2) | Out-Null

# 5) Add an empty "_GoBack" bookmark right at the end of the (now edited) first
#    paragraph's text, i.e. right after "...installazione:" and before the
#    paragraph mark. A genuinely-collapsed range sitting exactly on that
#    paragraph-end boundary mis-resolves, so we briefly insert a one-character
#    placeholder, bookmark the (non-collapsed) range around it, then clear the
#    placeholder text back out -- this leaves the bookmark start/end together
#    at the correct spot.
$installParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "Il programma non necessita di alcuna installazione:*") {
        $installParaIndex = $i
        break
    }
}
$p = $d.Paragraphs.Item($installParaIndex)
$end = $p.Range
$end.Collapse(0)
$end.Move(1, -1)
$end.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $end) | Out-Null
$end.Text = ""
